$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-01 Wednesday", "2023-11-02 Thursday"),
    @("89×92=", "85×98="),
    @("77×47=", "68×48="),
    @("57×32=", "43×80="),
    @("60×59=", "92×17="),
    @("12×27=", "84×71="),
    @("79×63=", "90×56="),
    @("42×48=", "68×28="),
    @("19×40=", "95×31="),
    @("76×47=", "49×62="),
    @("97×97=", "94×92="),
    @("83×94=", "14×79="),
    @("16×15=", "56×78="),
    @("88×68=", "31×11="),
    @("89×98=", "85×46="),
    @("88×95=", "11×95="),
    @("78×84=", "67×26="),
    @("45×43=", "30×69="),
    @("70×49=", "21×19="),
    @("79×95=", "93×34="),
    @("94×45=", "37×46="),
    @("39×67=", "99×54="),
    @("70×34=", "89×97="),
    @("22×83=", "67×81="),
    @("38×37=", "25×40="),
    @("64×16=", "50×53=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
